$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 130419.375
$ws.Range("I80").Value = 421.5
$ws.Range("K80").Value = 1264.5
$ws.Range("M80").Value = -266.5
$ws.Range("H83").Value = 130419.375
$ws.Range("I83").Value = 421.5
$ws.Range("K83").Value = 3793.5
$ws.Range("M83").Value = 1198.5
$ws.Range("H92").Value = 569.7143
$ws.Range("I92").Value = 540.25
$ws.Range("J92").Value = 609
$ws.Range("K92").Value = 540.25
$ws.Range("L92").Value = 609
$ws.Range("M92").Value = 707.75
$ws.Range("N92").Value = -3105
$ws.Range("H98").Value = 2977.1177
$ws.Range("I98").Value = 782.63635
$ws.Range("K98").Value = 782.63635
$ws.Range("M98").Value = 715.36365
$ws.Range("H111").Value = 3105.2
$ws.Range("I111").Value = 529
$ws.Range("K111").Value = 1587
$ws.Range("M111").Value = 1480
$ws.Range("H118").Value = 1870
$ws.Range("I118").Value = 1870
$ws.Range("K118").Value = 5610
$ws.Range("M118").Value = -3953
$ws.Range("H122").Value = 2977.1177
$ws.Range("I122").Value = 782.63635
$ws.Range("K122").Value = 2347.90905
$ws.Range("M122").Value = 102.0909499999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 5000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 5000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 5000
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -6040
$ws.Range("H61").Value = 3920.4443
$ws.Range("I61").Value = 1148.3125
$ws.Range("J61").Value = 7952.636
$ws.Range("K61").Value = 1148.3125
$ws.Range("L61").Value = 7952.636
$ws.Range("M61").Value = -936.3125
$ws.Range("N61").Value = -8376.636
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H112").Value = 36596.5
$ws.Range("J112").Value = 36596.5
$ws.Range("L112").Value = 36596.5
$ws.Range("N112").Value = -39550.5
$ws.Range("H134").Value = 79996.5
$ws.Range("J134").Value = 79996.5
$ws.Range("L134").Value = 79996.5
$ws.Range("N134").Value = -90136.5
$ws.Range("H136").Value = 3920.4443
$ws.Range("I136").Value = 1148.3125
$ws.Range("J136").Value = 7952.636
$ws.Range("K136").Value = 3444.9375
$ws.Range("L136").Value = 23857.908
$ws.Range("M136").Value = -894.9375
$ws.Range("N136").Value = -28957.908
$ws.Range("H139").Value = 39999.5
$ws.Range("J139").Value = 60000
$ws.Range("L139").Value = 60000
$ws.Range("N139").Value = -70280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 15000
$ws.Range("J38").Value = 15000
$ws.Range("L38").Value = 15000
$ws.Range("N38").Value = -15832
$ws.Range("H81").Value = 43693.75
$ws.Range("J81").Value = 43693.75
$ws.Range("L81").Value = 43693.75
$ws.Range("N81").Value = -45815.75
$ws.Range("H84").Value = 43693.75
$ws.Range("J84").Value = 43693.75
$ws.Range("L84").Value = 131081.25
$ws.Range("N84").Value = -141689.25
$ws.Range("H134").Value = 3327.111
$ws.Range("I134").Value = 3305.5
$ws.Range("K134").Value = 9916.5
$ws.Range("M134").Value = -7381.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2357.9285
$ws.Range("I16").Value = 1919
$ws.Range("J16").Value = 2796.8572
$ws.Range("K16").Value = 1919
$ws.Range("L16").Value = 2796.8572
$ws.Range("M16").Value = -1632
$ws.Range("N16").Value = -3370.8572
$ws.Range("H22").Value = 42909
$ws.Range("I22").Value = 1862.5
$ws.Range("K22").Value = 1862.5
$ws.Range("M22").Value = -1512.5
$ws.Range("H31").Value = 2056.9285
$ws.Range("I31").Value = 1874.8334
$ws.Range("K31").Value = 1874.8334
$ws.Range("M31").Value = -1579.8334
$ws.Range("H34").Value = 2056.9285
$ws.Range("I34").Value = 1874.8334
$ws.Range("K34").Value = 1874.8334
$ws.Range("M34").Value = -1672.8334
$ws.Range("H113").Value = 2357.9285
$ws.Range("I113").Value = 1919
$ws.Range("J113").Value = 2796.8572
$ws.Range("K113").Value = 1919
$ws.Range("L113").Value = 2796.8572
$ws.Range("M113").Value = 251
$ws.Range("N113").Value = -7136.8572
$ws.Range("H134").Value = 4384.4287
$ws.Range("I134").Value = 4281.8335
$ws.Range("K134").Value = 12845.5005
$ws.Range("M134").Value = -10310.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 776
$ws.Range("I18").Value = 776
$ws.Range("K18").Value = 2328
$ws.Range("M18").Value = -2159
$ws.Range("H113").Value = 1084.9412
$ws.Range("I113").Value = 1139.3334
$ws.Range("K113").Value = 3418.0002
$ws.Range("M113").Value = -1248.0002
$ws.Range("H124").Value = 1239.4
$ws.Range("I124").Value = 1239.4
$ws.Range("K124").Value = 3718.2
$ws.Range("M124").Value = 1191.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 118.933334
$ws.Range("I2").Value = 146.58333
$ws.Range("K2").Value = 146.58333
$ws.Range("M2").Value = -33.58332999999999
$ws.Range("H102").Value = 1261.7368
$ws.Range("I102").Value = 929.625
$ws.Range("J102").Value = 3033
$ws.Range("K102").Value = 929.625
$ws.Range("L102").Value = 3033
$ws.Range("M102").Value = 692.375
$ws.Range("N102").Value = -6277
$ws.Range("H126").Value = 2670.5
$ws.Range("I126").Value = 2670.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8011.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5541.5
$ws.Range("N126").ClearContents()
$ws.Range("H141").Value = 49999
$ws.Range("J141").Value = 49999
$ws.Range("L141").Value = 49999
$ws.Range("N141").Value = -60359

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6816.3213
$ws.Range("I7").Value = 3415.5454
$ws.Range("J7").Value = 9016.823
$ws.Range("K7").Value = 3415.5454
$ws.Range("L7").Value = 9016.823
$ws.Range("M7").Value = -3303.5454
$ws.Range("N7").Value = -9240.823
$ws.Range("H22").Value = 1509.8
$ws.Range("I22").Value = 899.5
$ws.Range("J22").Value = 1916.6666
$ws.Range("K22").Value = 899.5
$ws.Range("L22").Value = 1916.6666
$ws.Range("M22").Value = -604.5
$ws.Range("N22").Value = -2506.6666
$ws.Range("H27").Value = 1509.8
$ws.Range("I27").Value = 899.5
$ws.Range("J27").Value = 1916.6666
$ws.Range("K27").Value = 899.5
$ws.Range("L27").Value = 1916.6666
$ws.Range("M27").Value = -792.5
$ws.Range("N27").Value = -2130.6666
$ws.Range("H46").Value = 43200.25
$ws.Range("J46").Value = 1675.3
$ws.Range("L46").Value = 1675.3
$ws.Range("N46").Value = -2051.3
$ws.Range("H68").Value = 2667
$ws.Range("I68").Value = 1998.6666
$ws.Range("K68").Value = 1998.6666
$ws.Range("M68").Value = -1249.6666
$ws.Range("H71").Value = 2667
$ws.Range("I71").Value = 1998.6666
$ws.Range("K71").Value = 9993.333000000001
$ws.Range("M71").Value = -6249.333000000001
$ws.Range("H82").Value = 1433.6111
$ws.Range("J82").Value = 1585.625
$ws.Range("L82").Value = 1585.625
$ws.Range("N82").Value = -2307.625
$ws.Range("H85").Value = 1433.6111
$ws.Range("J85").Value = 1585.625
$ws.Range("L85").Value = 1585.625
$ws.Range("N85").Value = -4081.625
$ws.Range("H126").Value = 6816.3213
$ws.Range("I126").Value = 3415.5454
$ws.Range("J126").Value = 9016.823
$ws.Range("K126").Value = 10246.6362
$ws.Range("L126").Value = 27050.469
$ws.Range("M126").Value = -7776.636200000001
$ws.Range("N126").Value = -31990.469
$ws.Range("H136").Value = 3891.2778
$ws.Range("I136").Value = 3526.7693
$ws.Range("K136").Value = 10580.3079
$ws.Range("M136").Value = -8030.3079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H107").Value = 1221.25
$ws.Range("I107").Value = 486
$ws.Range("J107").Value = 1466.3334
$ws.Range("K107").Value = 1458
$ws.Range("L107").Value = 4399.0002
$ws.Range("M107").Value = 462
$ws.Range("N107").Value = -8239.0002
$ws.Range("H122").Value = 461.33334
$ws.Range("I122").Value = 447
$ws.Range("J122").Value = 490
$ws.Range("K122").Value = 1341
$ws.Range("L122").Value = 1470
$ws.Range("M122").Value = 1109
$ws.Range("N122").Value = -6370
$ws.Range("H136").Value = 1559.1818
$ws.Range("I136").Value = 1537.9524
$ws.Range("J136").Value = 2005
$ws.Range("K136").Value = 4613.857199999999
$ws.Range("L136").Value = 6015
$ws.Range("M136").Value = -2063.857199999999
$ws.Range("N136").Value = -11115

